# The worksheet holds a 7-column x 3-row grid of "keywordN" suggestions
# (rows 2-4, columns A-G), with day-of-week headers in row 1.
#
# Previously the keywords were written down each column before moving to
# the next column (A2=keyword1, A3=keyword2, A4=keyword3, B2=keyword4, ...).
# The script was updated to write the longest/shortest suggestions across
# each row instead, so the keywords now read left-to-right, top-to-bottom:
# row 2 = keyword1..keyword7, row 3 = keyword8..keyword14,
# row 4 = keyword15..keyword21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$keywords = @(
    "keyword1","keyword2","keyword3","keyword4","keyword5","keyword6","keyword7",
    "keyword8","keyword9","keyword10","keyword11","keyword12","keyword13","keyword14",
    "keyword15","keyword16","keyword17","keyword18","keyword19","keyword20","keyword21"
)

$columns = @("A","B","C","D","E","F","G")

$i = 0
for ($row = 2; $row -le 4; $row++) {
    foreach ($col in $columns) {
        $ws.Range("$col$row").Value = $keywords[$i]
        $i++
    }
}
